# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet right
# before the old "Late" column (column N), shifting Late/Date/Outstanding
# one column to the right, then make that sheet the active/selected one
# (matching the author's saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position 14 (N). Everything from the old
# N onward (Late / Date / Outstanding) shifts right to O / P / Q.
$ws.Columns.Item(14).EntireColumn.Insert() | Out-Null

# The newly inserted column keeps the width of its left neighbour (M),
# matching Excel's default "inherit width from the column to the left"
# behaviour when inserting columns.
$ws.Columns.Item(14).ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab and set its selection,
# matching the saved view state in the workbook.
$ws.Activate() | Out-Null
$ws.Range("K16").Select() | Out-Null
